$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update header text: volume number and report week dates ---
$ws.Range("A8").Value = "Volume 30   Number  28"
$ws.Range("C9").Value = "Report Covering the Week  7/10/2023  Through  7/16/2023"

# --- Update crime statistics table (rows 14-29) ---

$ws.Range("L14").Value = 50
$ws.Range("N14").Value = -84.210526315789
$ws.Range("C14").Copy()
$ws.Range("G15").PasteSpecial(-4122)
$ws.Range("G15").Value = "0"
$ws.Range("E14").Copy()
$ws.Range("H15").PasteSpecial(-4122)
$ws.Range("H15").Value = "***.*"
$ws.Range("N15").Value = -85
$ws.Range("C16").Value = 3
$ws.Range("E16").Value = 0
$ws.Range("F16").Value = 15
$ws.Range("G16").Value = 12
$ws.Range("H16").Value = 25
$ws.Range("I16").Value = 92
$ws.Range("J16").Value = 95
$ws.Range("K16").Value = -3.157894736842
$ws.Range("L16").Value = 12.195121951219
$ws.Range("M16").Value = -14.018691588785
$ws.Range("N16").Value = -79.325842696629
$ws.Range("C17").Value = 7
$ws.Range("D17").Value = 4
$ws.Range("E17").Value = 75
$ws.Range("F17").Value = 20
$ws.Range("G17").Value = 23
$ws.Range("H17").Value = -13.043478260869
$ws.Range("I17").Value = 133
$ws.Range("J17").Value = 124
$ws.Range("K17").Value = 7.258064516129
$ws.Range("L17").Value = -0.746268656716
$ws.Range("M17").Value = 20.909090909090
$ws.Range("N17").Value = -65.984654731457
$ws.Range("C18").Value = 1
$ws.Range("I14").Copy()
$ws.Range("D18").PasteSpecial(-4122)
$ws.Range("D18").Value = 4
$ws.Range("K14").Copy()
$ws.Range("E18").PasteSpecial(-4122)
$ws.Range("E18").Value = -75
$ws.Range("F18").Value = 6
$ws.Range("G18").Value = 7
$ws.Range("H18").Value = -14.285714285714
$ws.Range("I18").Value = 64
$ws.Range("J18").Value = 92
$ws.Range("K18").Value = -30.434782608695
$ws.Range("L18").Value = -22.891566265060
$ws.Range("M18").Value = 3.225806451612
$ws.Range("N18").Value = -86.721991701244
$ws.Range("C19").Value = 6
$ws.Range("E19").Value = -33.333333333333
$ws.Range("F19").Value = 31
$ws.Range("G19").Value = 39
$ws.Range("H19").Value = -20.512820512820
$ws.Range("I19").Value = 188
$ws.Range("J19").Value = 208
$ws.Range("K19").Value = -9.615384615384
$ws.Range("L19").Value = -2.083333333333
$ws.Range("M19").Value = 16.049382716049
$ws.Range("N19").Value = -16.814159292035
$ws.Range("C14").Copy()
$ws.Range("C20").PasteSpecial(-4122)
$ws.Range("C20").Value = "0"
$ws.Range("D20").Value = 2
$ws.Range("E20").Value = -100
$ws.Range("F20").Value = 9
$ws.Range("G20").Value = 6
$ws.Range("H20").Value = 50
$ws.Range("J20").Value = 29
$ws.Range("K20").Value = 24.137931034482
$ws.Range("L20").Value = 16.129032258064
$ws.Range("M20").Value = 176.923076923077
$ws.Range("N20").Value = -53.846153846153
$ws.Range("C21").Value = 17
$ws.Range("D21").Value = 22
$ws.Range("E21").Value = -22.727272727272
$ws.Range("F21").Value = 81
$ws.Range("G21").Value = 87
$ws.Range("H21").Value = -6.896551724137
$ws.Range("I21").Value = 519
$ws.Range("J21").Value = 554
$ws.Range("K21").Value = -6.317689530685
$ws.Range("L21").Value = -1.890359168241
$ws.Range("M21").Value = 11.612903225806
$ws.Range("N21").Value = -68.753762793497
$ws.Range("C14").Copy()
$ws.Range("C22").PasteSpecial(-4122)
$ws.Range("C22").Value = "0"
$ws.Range("C14").Copy()
$ws.Range("D22").PasteSpecial(-4122)
$ws.Range("D22").Value = "0"
$ws.Range("E14").Copy()
$ws.Range("E22").PasteSpecial(-4122)
$ws.Range("E22").Value = "***.*"
$ws.Range("F22").Value = 1
$ws.Range("G22").Value = 1
$ws.Range("I22").Value = 10
$ws.Range("K22").Value = -52.380952380952
$ws.Range("L22").Value = -33.333333333333
$ws.Range("M22").Value = 25
$ws.Range("I14").Copy()
$ws.Range("C23").PasteSpecial(-4122)
$ws.Range("C23").Value = 2
$ws.Range("D23").Value = 1
$ws.Range("E23").Value = 100
$ws.Range("F23").Value = 6
$ws.Range("G23").Value = 9
$ws.Range("H23").Value = -33.333333333333
$ws.Range("I23").Value = 42
$ws.Range("J23").Value = 55
$ws.Range("K23").Value = -23.636363636363
$ws.Range("L23").Value = -2.325581395348
$ws.Range("C24").Value = 23
$ws.Range("E24").Value = -20.689655172413
$ws.Range("F24").Value = 132
$ws.Range("G24").Value = 100
$ws.Range("H24").Value = 32
$ws.Range("I24").Value = 711
$ws.Range("J24").Value = 621
$ws.Range("K24").Value = 14.492753623188
$ws.Range("L24").Value = 13.036565977742
$ws.Range("M24").Value = 38.596491228070
$ws.Range("C25").Value = 6
$ws.Range("D25").Value = 8
$ws.Range("F25").Value = 33
$ws.Range("G25").Value = 37
$ws.Range("H25").Value = -10.810810810810
$ws.Range("I25").Value = 230
$ws.Range("J25").Value = 237
$ws.Range("K25").Value = -2.953586497890
$ws.Range("L25").Value = 4.545454545454
$ws.Range("M25").Value = -20.689655172413
$ws.Range("I14").Copy()
$ws.Range("D26").PasteSpecial(-4122)
$ws.Range("D26").Value = 1
$ws.Range("K14").Copy()
$ws.Range("E26").PasteSpecial(-4122)
$ws.Range("E26").Value = -100
$ws.Range("G26").Value = 1
$ws.Range("H26").Value = 0
$ws.Range("J26").Value = 9
$ws.Range("K26").Value = -11.111111111111
$ws.Range("C27").Value = 1
$ws.Range("I14").Copy()
$ws.Range("D27").PasteSpecial(-4122)
$ws.Range("D27").Value = 1
$ws.Range("K14").Copy()
$ws.Range("E27").PasteSpecial(-4122)
$ws.Range("E27").Value = 0
$ws.Range("G27").Value = 2
$ws.Range("H27").Value = 50
$ws.Range("I27").Value = 22
$ws.Range("J27").Value = 29
$ws.Range("K27").Value = -24.137931034482
$ws.Range("L27").Value = 0
$ws.Range("I14").Copy()
$ws.Range("C28").PasteSpecial(-4122)
$ws.Range("C28").Value = 1
$ws.Range("C14").Copy()
$ws.Range("D28").PasteSpecial(-4122)
$ws.Range("D28").Value = "0"
$ws.Range("E14").Copy()
$ws.Range("E28").PasteSpecial(-4122)
$ws.Range("E28").Value = "***.*"
$ws.Range("I14").Copy()
$ws.Range("F28").PasteSpecial(-4122)
$ws.Range("F28").Value = 1
$ws.Range("H28").Value = -66.666666666666
$ws.Range("I28").Value = 7
$ws.Range("K28").Value = -30
$ws.Range("L28").Value = -56.25
$ws.Range("M28").Value = -50
$ws.Range("N28").Value = -88.135593220339
$ws.Range("I14").Copy()
$ws.Range("C29").PasteSpecial(-4122)
$ws.Range("C29").Value = 1
$ws.Range("C14").Copy()
$ws.Range("D29").PasteSpecial(-4122)
$ws.Range("D29").Value = "0"
$ws.Range("E14").Copy()
$ws.Range("E29").PasteSpecial(-4122)
$ws.Range("E29").Value = "***.*"
$ws.Range("I14").Copy()
$ws.Range("F29").PasteSpecial(-4122)
$ws.Range("F29").Value = 1
$ws.Range("H29").Value = -66.666666666666
$ws.Range("I29").Value = 6
$ws.Range("K29").Value = -40
$ws.Range("L29").Value = -62.5
$ws.Range("M29").Value = -50
$ws.Range("N29").Value = -88.888888888888
